$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(358, "Thu Oct 17 01:59:50 2024", 118, 128),
    @(359, "Thu Oct 17 01:59:53 2024", 118, 128),
    @(360, "Thu Oct 17 01:59:56 2024", 118, 128),
    @(361, "Thu Oct 17 01:59:59 2024", 118, 128),
    @(362, "Thu Oct 17 02:00:02 2024", 118, 129),
    @(363, "Thu Oct 17 02:00:05 2024", 119, 129),
    @(364, "Thu Oct 17 02:05:56 2024", 116, 127),
    @(365, "Thu Oct 17 02:05:59 2024", 115, 126),
    @(366, "Thu Oct 17 02:06:02 2024", 115, 127)
)

$startRow = 359
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rec = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
}
